$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 5285.2856
$ws.Range("I26").Value = 4499.5
$ws.Range("K26").Value = 4499.5
$ws.Range("M26").Value = -4169.5

$ws.Range("H39").Value = 7508.5
$ws.Range("I39").Value = 5016
$ws.Range("K39").Value = 5016
$ws.Range("M39").Value = -4496

$ws.Range("H61").Value = 58942930
$ws.Range("I61").Value = 83418240
$ws.Range("J61").Value = 202176
$ws.Range("K61").Value = 83418240
$ws.Range("L61").Value = 202176
$ws.Range("M61").Value = -83418028
$ws.Range("N61").Value = -202600

$ws.Range("H74").Value = 19385606
$ws.Range("I74").Value = 41834190
$ws.Range("J74").Value = 143959.14
$ws.Range("K74").Value = 41834190
$ws.Range("L74").Value = 143959.14
$ws.Range("M74").Value = -41833316
$ws.Range("N74").Value = -145707.14

$ws.Range("H77").Value = 19385606
$ws.Range("I77").Value = 41834190
$ws.Range("J77").Value = 143959.14
$ws.Range("K77").Value = 209170950
$ws.Range("L77").Value = 719795.7000000001
$ws.Range("M77").Value = -209166582
$ws.Range("N77").Value = -728531.7000000001

$ws.Range("H122").Value = 1218.7778
$ws.Range("I122").Value = 1076.5
$ws.Range("J122").Value = 2357
$ws.Range("K122").Value = 3229.5
$ws.Range("L122").Value = 7071
$ws.Range("M122").Value = -779.5
$ws.Range("N122").Value = -11971

$ws.Range("H136").Value = 58942930
$ws.Range("I136").Value = 83418240
$ws.Range("J136").Value = 202176
$ws.Range("K136").Value = 250254720
$ws.Range("L136").Value = 606528
$ws.Range("M136").Value = -250252170
$ws.Range("N136").Value = -611628

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2800000
$ws.Range("I7").Value = 2800000
$ws.Range("K7").Value = 2800000
$ws.Range("M7").Value = -2799887

$ws.Range("H134").Value = 3131.4285
$ws.Range("I134").Value = 10000
$ws.Range("J134").Value = 1986.6666
$ws.Range("K134").Value = 30000
$ws.Range("L134").Value = 5959.9998
$ws.Range("M134").Value = -27465
$ws.Range("N134").Value = -11029.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4112.304
$ws.Range("I86").Value = 4098.706
$ws.Range("J86").Value = 4150.8335
$ws.Range("K86").Value = 4098.706
$ws.Range("L86").Value = 4150.8335
$ws.Range("M86").Value = -2975.706
$ws.Range("N86").Value = -6396.8335

$ws.Range("H89").Value = 4112.304
$ws.Range("I89").Value = 4098.706
$ws.Range("J89").Value = 4150.8335
$ws.Range("K89").Value = 20493.53
$ws.Range("L89").Value = 20754.1675
$ws.Range("M89").Value = -14877.53
$ws.Range("N89").Value = -31986.1675

$ws.Range("H134").Value = 70105.125
$ws.Range("I134").Value = 3345.6365
$ws.Range("J134").Value = 216976
$ws.Range("K134").Value = 10036.9095
$ws.Range("L134").Value = 650928
$ws.Range("M134").Value = -7501.9095
$ws.Range("N134").Value = -655998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4250
$ws.Range("J62").Value = 4250
$ws.Range("L62").Value = 12750
$ws.Range("N62").Value = -14122

$ws.Range("H65").Value = 4250
$ws.Range("J65").Value = 4250
$ws.Range("L65").Value = 38250
$ws.Range("N65").Value = -45114

$ws.Range("H68").Value = 1137.5
$ws.Range("I68").Value = 731.4815
$ws.Range("J68").Value = 2355.5557
$ws.Range("K68").Value = 2194.4445
$ws.Range("L68").Value = 7066.6671
$ws.Range("M68").Value = -1383.4445
$ws.Range("N68").Value = -8688.667099999999

$ws.Range("H71").Value = 1137.5
$ws.Range("I71").Value = 731.4815
$ws.Range("J71").Value = 2355.5557
$ws.Range("K71").Value = 6583.3335
$ws.Range("L71").Value = 21200.0013
$ws.Range("M71").Value = -2527.3335
$ws.Range("N71").Value = -29312.0013

$ws.Range("H123").Value = 975
$ws.Range("I123").Value = 975
$ws.Range("K123").Value = 2925
$ws.Range("M123").Value = -475

$ws.Range("H138").Value = 2934.0588
$ws.Range("I138").Value = 2164.3333
$ws.Range("K138").Value = 6492.999899999999
$ws.Range("M138").Value = -1352.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2112.6667
$ws.Range("I102").Value = 2122.9092
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2122.9092
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -500.9092000000001
$ws.Range("N102").Value = -5244

$ws.Range("H132").Value = 71805.92999999999
$ws.Range("I132").Value = 60970
$ws.Range("J132").Value = 87156.836
$ws.Range("K132").Value = 182910
$ws.Range("L132").Value = 261470.508
$ws.Range("M132").Value = -180380
$ws.Range("N132").Value = -266530.508

$ws.Range("H138").Value = 38500
$ws.Range("J138").Value = 38500
$ws.Range("L138").Value = 38500
$ws.Range("N138").Value = -48780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 869.3
$ws.Range("I9").Value = 218.57143
$ws.Range("J9").Value = 2387.6667
$ws.Range("K9").Value = 218.57143
$ws.Range("L9").Value = 2387.6667
$ws.Range("M9").Value = 5.428570000000008
$ws.Range("N9").Value = -2835.6667

$ws.Range("H64").Value = 30263.666
$ws.Range("J64").Value = 30263.666
$ws.Range("L64").Value = 30263.666
$ws.Range("N64").Value = -30713.666

$ws.Range("H67").Value = 30263.666
$ws.Range("J67").Value = 30263.666
$ws.Range("L67").Value = 30263.666
$ws.Range("N67").Value = -31823.666

$ws.Range("H70").Value = 39442
$ws.Range("J70").Value = 39442
$ws.Range("L70").Value = 39442
$ws.Range("N70").Value = -39982

$ws.Range("H73").Value = 39442
$ws.Range("J73").Value = 39442
$ws.Range("L73").Value = 39442
$ws.Range("N73").Value = -41314

$ws.Range("H122").Value = 3201.5264
$ws.Range("I122").Value = 3266.2856
$ws.Range("J122").Value = 3163.75
$ws.Range("K122").Value = 9798.856800000001
$ws.Range("L122").Value = 9491.25
$ws.Range("M122").Value = -7348.856800000001
$ws.Range("N122").Value = -14391.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

$ws.Range("H122").Value = 2236.9062
$ws.Range("I122").Value = 1695.0625
$ws.Range("J122").Value = 2778.75
$ws.Range("K122").Value = 5085.1875
$ws.Range("L122").Value = 8336.25
$ws.Range("M122").Value = -2635.1875
$ws.Range("N122").Value = -13236.25

$ws.Range("H136").Value = 75534.19
$ws.Range("I136").Value = 46705.59
$ws.Range("J136").Value = 202380
$ws.Range("K136").Value = 140116.77
$ws.Range("L136").Value = 607140
$ws.Range("M136").Value = -137566.77
$ws.Range("N136").Value = -612240
